$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.408.63"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.69%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.723.13"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.43%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "243.50"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.08%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4919"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +2.14%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2610"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -2.45%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06203"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.34%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.720.15"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -0.64%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.06987"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -1.95%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.44"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -1.27%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.539"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.23%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5995"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -2.12%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "77.43"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.07%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.000"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.06%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.390.90"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.78%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.000"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.08%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007205"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +3.37%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.34"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -2.03%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.942.97"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.55%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.470"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -1.20%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -2.72%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.154"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -1.77%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "137.81"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.38%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.26"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.71%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -1.26%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "106.90"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -1.84%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -3.47%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.948"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.80%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08012"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.07%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.677"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.33%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04515"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.53%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9994"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.02%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.47%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9989"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.65%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.6260"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -1.40%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.9467"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +5.30%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.393"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.66%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.948"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -5.35%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9998"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.36%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.01481"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -1.51%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "99.67"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -3.37%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.287"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -3.33%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.3849"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -1.45%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "6.815"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -5.65%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1168"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.41%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.762"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -2.11%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "30.20"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -1.63%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.235"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -1.63%  "
